# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 19:22"

# Row 4 - Estados Unidos: refreshed totals
$ws.Range("B4").Value = 1075582
$ws.Range("C4").Value = 11388
$ws.Range("D4").Value = 149686
$ws.Range("E4").Value = 863581
$ws.Range("G4").Value = 660
$ws.Range("H4").Value = 62315

# Reino Unido overtakes Francia in the ranking: row 7 now holds Reino Unido's
# refreshed figures, row 8 now holds Francia's (previously row 7's) figures.
$ws.Range("A7").Value = "Reino Unido"
$ws.Range("B7").Value = 171253
$ws.Range("C7").Value = 6032
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 144138
$ws.Range("F7").Value = 1559
$ws.Range("G7").Value = 674
$ws.Range("H7").Value = 26771

$ws.Range("A8").Value = "Francia"
$ws.Range("B8").Value = 166420
$ws.Range("D8").Value = 48228
$ws.Range("E8").Value = 94105
$ws.Range("F8").Value = 4207
$ws.Range("H8").Value = 24087

# India overtakes Peru in the ranking: row 18 now holds India's refreshed
# figures, row 19 now holds Peru's (previously row 18's) figures.
$ws.Range("A18").Value = "India"
$ws.Range("B18").Value = 34780
$ws.Range("C18").Value = 1718
$ws.Range("D18").Value = 9068
$ws.Range("E18").Value = 24561
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 72
$ws.Range("H18").Value = 1151

$ws.Range("A19").Value = "Peru"
$ws.Range("B19").Value = 33931
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 10037
$ws.Range("E19").Value = 22951
$ws.Range("F19").Value = 623
$ws.Range("H19").Value = 943

# Row 22 - Ecuador: refreshed totals
$ws.Range("B22").Value = 24934
$ws.Range("C22").Value = 259
$ws.Range("D22").Value = 1558
$ws.Range("E22").Value = 22476
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = 900

# Row 25 - Irlanda: refreshed totals
$ws.Range("B25").Value = 20612
$ws.Range("C25").Value = 359
$ws.Range("E25").Value = 5994
$ws.Range("F25").Value = 123
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 1232

# Row 109 - Georgia: refreshed totals
$ws.Range("D109").Value = 184
$ws.Range("E109").Value = 349
